$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value()
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 3.34 = 12760.54 pesos`n✅ 12760.54 pesos = 3.32 = 973.91 Bs"), "✅ 1000 Bs = 3.34 = 12709.03 pesos`n✅ 12709.03 pesos = 3.34 = 953.09 Bs"
$cell.Value = $text

# --- tasas: update rate tracking values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("O10").Value = 3800
$wsTasas.Range("N12").Value = 3807
$wsTasas.Range("O12").Value = 285.5
